$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header details ---
$ws.Range("C2").Value = "Hartmut"

# B3 holds a long digit string that must stay text (not become a number).
# Stage it as text in a scratch cell (forced to Text number format) and
# paste-special as values so the numeric-looking string is not re-typed
# (which would make Excel infer it as a numeric value).
$scratch = $ws.Range("Z1")
$scratch.NumberFormat = "@"
$scratch.Value = "2570314725427075"
$scratch.Copy()
$ws.Range("B3").PasteSpecial(-4163)  # xlPasteValues
$scratch.Clear()

$ws.Range("C3").Value = "Mohaupt"

# --- Opening balance line ---
$ws.Range("D5").Value = "KONTOSTAND AM 15.09.2024"

# --- Row 6 ---
$ws.Range("B6").Value = "17.09."
$ws.Range("C6").Value = "18.09."
$ws.Range("D6").Value = "BEITRAG Allianz SE K-89310533"
$ws.Range("E6").Value = "56,04-"

# --- Row 7 ---
$ws.Range("B7").Value = "21.09."
$ws.Range("C7").Value = "22.09."
$ws.Range("D7").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 68057210"
$ws.Range("E7").Value = "86,28-"

# --- Row 8 ---
$ws.Range("B8").Value = "25.09."
$ws.Range("C8").Value = "26.09."
$ws.Range("D8").Value = "KARTENZ./25.09 ALDI SUED RO"
$ws.Range("E8").Value = "127,55-"

# --- Row 9 ---
$ws.Range("B9").Value = "26.09."
$ws.Range("C9").Value = "27.09."
$ws.Range("D9").Value = "RECHNUNG VODAFONE GMBH 71983496"
$ws.Range("E9").Value = "40,36-"

# --- Row 10 ---
$ws.Range("B10").Value = "27.09."
$ws.Range("C10").Value = "28.09."
$ws.Range("D10").Value = "KARTENZAHLUNG SHELL TANKSTELLE"
$ws.Range("E10").Value = "73,33-"

# --- Row 11: transaction removed, row becomes blank ---
# B11/C11/D11 lose their values (style stays s=8, same as B12/C12 which are
# the already-blank counterparts in this table).
$ws.Range("B11").Value = ""
$ws.Range("C11").Value = ""
$ws.Range("D11").Value = ""

# E11 becomes blank too, but its style changes from the "right aligned only"
# style (shared with E6:E10) to the "right aligned + vertically centered +
# wrap" style already used elsewhere in the sheet (e.g. the BESCHREIBUNG
# cells in column D). Grab that exact format via copy/paste-special so we
# land on the existing style instead of synthesizing a brand-new one, then
# nudge horizontal alignment back to right (it is the only attribute that
# differs) to reach the target style precisely.
$ws.Range("D6").Copy()
$ws.Range("E11").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("E11").HorizontalAlignment = -4152  # xlRight
$ws.Range("E11").Value = ""

# --- Closing balance line ---
$ws.Range("D12").Value = "KONTOSTAND AM 02.10.2024"
$ws.Range("E12").Value = "383,56-"

# --- Next statement date ---
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 10.10.2024"

$excel.CutCopyMode = $false
